$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Row 2
$ws.Range("D2").Value = "43.610.14"
$ws.Range("E2").Value = "  +0.03%  "

# Row 3
$ws.Range("D3").Value = "2.316.59"

# Row 4
$ws.Range("E4").Value = "  +0.07%  "

# Row 5
$ws.Range("D5").Value = "'269.27"
$ws.Range("E5").Value = "  -0.40%  "

# Row 6
$ws.Range("D6").Value = "'94.16"
$ws.Range("E6").Value = "  +5.08%  "

# Row 7
$ws.Range("E7").Value = "  +0.45%  "

# Row 8
$ws.Range("E8").Value = "  +0.02%  "

# Row 9
$ws.Range("D9").Value = "'0.618"
$ws.Range("E9").Value = "  +2.01%  "

# Row 10
$ws.Range("D10").Value = "'44.93"
$ws.Range("E10").Value = "  -1.73%  "

# Row 11
$ws.Range("E11").Value = "  +2.79%  "

# Row 12
$ws.Range("D12").Value = "'8.05"
$ws.Range("E12").Value = "  +2.28%  "

# Row 13
$ws.Range("E13").Value = "  +0.30%  "

# Row 14
$ws.Range("D14").Value = "2.663.15"
$ws.Range("E14").Value = "  +4.08%  "

# Row 15
$ws.Range("D15").Value = "'15.60"
$ws.Range("E15").Value = "  +3.46%  "

# Row 16
$ws.Range("D16").Value = "'0.856"
$ws.Range("E16").Value = "  +8.37%  "

# Row 17
$ws.Range("D17").Value = "2.321.55"
$ws.Range("E17").Value = "  +5.25%  "

# Row 18
$ws.Range("D18").Value = "43.519.71"
$ws.Range("E18").Value = "  -0.13%  "

# Row 19
$ws.Range("E19").Value = "  +6.46%  "

# Row 20
$ws.Range("D20").Value = "'6.35"
$ws.Range("E20").Value = "  +6.47%  "

# Row 21
$ws.Range("D21").Value = "'71.99"
$ws.Range("E21").Value = "  +2.38%  "

# Row 22
$ws.Range("D22").Value = "'240.18"
$ws.Range("E22").Value = "  +4.33%  "

# Row 23
$ws.Range("D23").Value = "'2.26"
$ws.Range("E23").Value = "  -3.43%  "

# Row 24
$ws.Range("D24").Value = "'9.36"
$ws.Range("E24").Value = "  +8.56%  "

# Row 25
$ws.Range("E25").Value = "  -0.07%  "

# Row 26
$ws.Range("B26").Value = "Cosmos"
$ws.Range("C26").Value = "https://coinranking.com/coin/Knsels4_Ol-Ny+cosmos-atom"
$ws.Range("D26").Value = "'11.41"
$ws.Range("E26").Value = "  +4.77%  "

# Row 27
$ws.Range("B27").Value = "PancakeSwap"
$ws.Range("C27").Value = "https://coinranking.com/coin/ncYFcP709+pancakeswap-cake"
$ws.Range("D27").Value = "'2.52"
$ws.Range("E27").Value = "  +1.25%  "

# Row 28
$ws.Range("D28").Value = "'3.48"
$ws.Range("E28").Value = "  -2.27%  "

# Row 29
$ws.Range("E29").Value = "  +0.62%  "

# Row 30
$ws.Range("D30").Value = "'38.16"
$ws.Range("E30").Value = "  -1.24%  "

# Row 31
$ws.Range("D31").Value = "'22.48"
$ws.Range("E31").Value = "  +8.64%  "

# Row 32
$ws.Range("D32").Value = "'172.23"
$ws.Range("E32").Value = "  -0.08%  "

# Row 33
$ws.Range("D33").Value = "'0.0895"
$ws.Range("E33").Value = "  -1.36%  "

# Row 34
$ws.Range("D34").Value = "'5.45"
$ws.Range("E34").Value = "  +1.95%  "

# Row 35
$ws.Range("E35").Value = "  +2.43%  "

# Row 36
$ws.Range("D36").Value = "'0.0358"
$ws.Range("E36").Value = "  +2.46%  "

# Row 37
$ws.Range("E37").Value = "  -3.11%  "

# Row 38
$ws.Range("D38").Value = "'4.37"
$ws.Range("E38").Value = "  +2.66%  "

# Row 39
$ws.Range("E39").Value = "  -1.40%  "

# Row 40
$ws.Range("E40").Value = "  +8.60%  "

# Row 41
$ws.Range("E41").Value = "  +10.63%  "

# Row 42
$ws.Range("D42").Value = "'1.35"
$ws.Range("E42").Value = "  +18.08%  "

# Row 43
$ws.Range("D43").Value = "'12.04"
$ws.Range("E43").Value = "  -3.45%  "

# Row 44
$ws.Range("D44").Value = "'9.17"
$ws.Range("E44").Value = "  +7.81%  "

# Row 45
$ws.Range("D45").Value = "'61.78"
$ws.Range("E45").Value = "  -1.98%  "

# Row 46
$ws.Range("D46").Value = "'5.33"
$ws.Range("E46").Value = "  -0.11%  "

# Row 47
$ws.Range("D47").Value = "'0.103"
$ws.Range("E47").Value = "  +4.17%  "

# Row 48
$ws.Range("D48").Value = "'100.16"
$ws.Range("E48").Value = "  +0.12%  "

# Row 49
$ws.Range("E49").Value = "  +2.77%  "

# Row 50
$ws.Range("D50").Value = "2.541.41"
$ws.Range("E50").Value = "  +3.97%  "

# Row 51
$ws.Range("D51").Value = "'0.185"
$ws.Range("E51").Value = "  +15.28%  "
